$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crosstab values (row 2 = Checked, row 3 = Unchecked)
# to match the reindexed/re-ordered output.
$ws.Range("B2").Value = 0.1901840490797546
$ws.Range("C2").Value = 0.2515337423312883
$ws.Range("D2").Value = 0.1104294478527607

$ws.Range("B3").Value = 0.2147239263803681
$ws.Range("C3").Value = 0.1717791411042945
$ws.Range("D3").Value = 0.06134969325153374
